# Update study-design / conflict-of-interest text for a handful of rows in
# the "PECO Full" sheet. The old placeholder values ("Both" and
# "Monte-Carlo") are replaced with the fuller descriptions used in the new
# graphs ("Both e-cigarrette and pharmaceutical" and "Simulation").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PECO Full")

# Conflict of Interest column (H) cells that said "Both" -> now fully spelled out
$bothCells = @("H3", "H4", "H9", "H12", "H28", "H33", "H42", "H95")
foreach ($cellRef in $bothCells) {
    $ws.Range($cellRef).Value = "Both e-cigarrette and pharmaceutical"
}

# Study design column (G) cell that said "Monte-Carlo" -> now "Simulation"
$ws.Range("G55").Value = "Simulation"

# Restore the view: scrolled to A49 with G53 selected (was D45 / G54)
$ws.Range("G53").Select()
